# Update "想去人数" (number of people interested) values for a few events
# that are listed in both the "展览" (Exhibitions) sheet and the
# "全部类型" (All types) aggregate sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 3-6, column F) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 308
$wsExpo.Range("F4").Value = 2820
$wsExpo.Range("F5").Value = 65
$wsExpo.Range("F6").Value = 597

# --- Sheet "全部类型" (rows 5-8, column F) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 308
$wsAll.Range("F6").Value = 2820
$wsAll.Range("F7").Value = 65
$wsAll.Range("F8").Value = 597
